$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7078968762894533
$ws.Range("C2").Value = 0.1948319280674866
$ws.Range("E2").Value = 0.4616680946494185
$ws.Range("F2").Value = 2.083719068667207
$ws.Range("G2").Value = 0.562933965054583
$ws.Range("H2").Value = 0.6548121233266926
$ws.Range("I2").Value = 0.4718468106774516
$ws.Range("J2").Value = 0.03454636694597069
$ws.Range("M2").Value = 0.5287727573242975
$ws.Range("N2").Value = 1.100234547277548
$ws.Range("B3").Value = 0.6239939310265186
$ws.Range("C3").Value = 0.1697473505300877
$ws.Range("E3").Value = 0.4593142221044602
$ws.Range("F3").Value = 2.064027292413527
$ws.Range("G3").Value = 0.5523729490152647
$ws.Range("H3").Value = 0.6550972489497866
$ws.Range("I3").Value = 0.4744858103288117
$ws.Range("J3").Value = 0.03513548234085517
$ws.Range("M3").Value = 0.4991156385457813
$ws.Range("N3").Value = 1.114557498325748
$ws.Range("B4").Value = 0.5724980932049277
$ws.Range("C4").Value = 0.1543125334117974
$ws.Range("E4").Value = 0.4580330418800642
$ws.Range("F4").Value = 2.053207255336829
$ws.Range("G4").Value = 0.5464143512934356
$ws.Range("H4").Value = 0.6556995105600407
$ws.Range("I4").Value = 0.4764763834131429
$ws.Range("J4").Value = 0.03551774972382393
$ws.Range("M4").Value = 0.4811173420663692
$ws.Range("N4").Value = 1.12384518782391
$ws.Range("B5").Value = 0.5515188867303209
$ws.Range("C5").Value = 0.1480144567486832
$ws.Range("E5").Value = 0.4575522895644504
$ws.Range("F5").Value = 2.049117104669847
$ws.Range("G5").Value = 0.5441177429705988
$ws.Range("H5").Value = 0.6560521571000066
$ws.Range("I5").Value = 0.4773803667398759
$ws.Range("J5").Value = 0.03567868378512173
$ws.Range("M5").Value = 0.4738362143664219
$ws.Range("N5").Value = 1.127753973411153
$ws.Range("B6").Value = 0.548035668843653
$ws.Range("C6").Value = 0.1469681663381266
$ws.Range("E6").Value = 0.4574749599472021
$ws.Range("F6").Value = 2.048457198262838
$ws.Range("G6").Value = 0.5437443195986731
$ws.Range("H6").Value = 0.6561171837615092
$ws.Range("I6").Value = 0.4775360695581909
$ws.Range("J6").Value = 0.03570571794238653
$ws.Range("M6").Value = 0.472630415061829
$ws.Range("N6").Value = 1.128410509986821
$ws.Range("B7").Value = 0.5722151362670616
$ws.Range("C7").Value = 0.1542276288275275
$ws.Range("E7").Value = 0.4580263908089606
$ws.Range("F7").Value = 2.053150802700358
$ws.Range("G7").Value = 0.5463828466102001
$ws.Range("H7").Value = 0.6557038326145346
$ws.Range("I7").Value = 0.4764881994501522
$ws.Range("J7").Value = 0.03551989926950538
$ws.Range("M7").Value = 0.4810189300670302
$ws.Range("N7").Value = 1.123897401200917
$ws.Range("B8").Value = 0.6789630780717744
$ws.Range("C8").Value = 0.1861895507846043
$ws.Range("E8").Value = 0.4608224609782354
$ws.Range("F8").Value = 2.076665307762568
$ws.Range("G8").Value = 0.5591829626984719
$ws.Range("H8").Value = 0.6548216473436526
$ws.Range("I8").Value = 0.4726797266591767
$ws.Range("J8").Value = 0.03474522320259954
$ws.Range("M8").Value = 0.5185032333260438
$ws.Range("N8").Value = 1.105070661741674
$ws.Range("B9").Value = 0.8884497822485855
$ws.Range("C9").Value = 0.248610651525297
$ws.Range("E9").Value = 0.4676051683522573
$ws.Range("F9").Value = 2.132886102750703
$ws.Range("G9").Value = 0.588489437321627
$ws.Range("H9").Value = 0.6564912915756054
$ws.Range("I9").Value = 0.4681621083632592
$ws.Range("J9").Value = 0.03338955243337915
$ws.Range("M9").Value = 0.5936822029872246
$ws.Range("N9").Value = 1.072069060816144
$ws.Range("B10").Value = 1.042453988174486
$ws.Range("C10").Value = 0.2943253505405323
$ws.Range("E10").Value = 0.473378412918251
$ws.Range("F10").Value = 2.180397772793668
$ws.Range("G10").Value = 0.6126329785488736
$ws.Range("H10").Value = 0.659805704984592
$ws.Range("I10").Value = 0.4666608708035014
$ws.Range("J10").Value = 0.03249374936624871
$ws.Range("M10").Value = 0.6499364827538869
$ws.Range("N10").Value = 1.050214902129827
$ws.Range("B11").Value = 1.112536092403559
$ws.Range("C11").Value = 0.3150928848662602
$ws.Range("E11").Value = 0.4761759448825629
$ws.Range("F11").Value = 2.203369647812238
$ws.Range("G11").Value = 0.6241942915337262
$ws.Range("H11").Value = 0.6617702433513841
$ws.Range("I11").Value = 0.4663767915339676
$ws.Range("J11").Value = 0.03210810447428303
$ws.Range("M11").Value = 0.6757501602928642
$ws.Range("N11").Value = 1.040793214191964
$ws.Range("B12").Value = 1.139077670180029
$ws.Range("C12").Value = 0.3229530379656183
$ws.Range("E12").Value = 0.4772598701401662
$ws.Range("F12").Value = 2.212264519284901
$ws.Range("G12").Value = 0.6286562143237404
$ws.Range("H12").Value = 0.6625801044551594
$ws.Range("I12").Value = 0.466326900573435
$ws.Range("J12").Value = 0.03196522549586422
$ws.Range("M12").Value = 0.685557153759774
$ws.Range("N12").Value = 1.037300334876534
$ws.Range("B13").Value = 1.133361341728175
$ws.Range("C13").Value = 0.3212603938201539
$ws.Range("E13").Value = 0.4770253361367409
$ws.Range("F13").Value = 2.210340126966557
$ws.Range("G13").Value = 0.627691515836247
$ws.Range("H13").Value = 0.6624027500406697
$ws.Range("I13").Value = 0.4663350747057464
$ws.Range("J13").Value = 0.03199585645342662
$ws.Range("M13").Value = 0.6834436260145509
$ws.Range("N13").Value = 1.038049254843521
$ws.Range("B14").Value = 1.114719626941962
$ws.Range("C14").Value = 0.3157396257266214
$ws.Range("E14").Value = 0.4762646282069696
$ws.Range("F14").Value = 2.204097504473282
$ws.Range("G14").Value = 0.6245596907542392
$ws.Range("H14").Value = 0.6618355482971197
$ws.Range("I14").Value = 0.4663715290086188
$ws.Range("J14").Value = 0.03209628640221185
$ws.Range("M14").Value = 0.6765563490778419
$ws.Range("N14").Value = 1.040504350724721
$ws.Range("B15").Value = 1.103301413214695
$ws.Range("C15").Value = 0.3123574672617337
$ws.Range("E15").Value = 0.4758018692215487
$ws.Range("F15").Value = 2.200299249316657
$ws.Range("G15").Value = 0.6226523064857616
$ws.Range("H15").Value = 0.6614967146270487
$ws.Range("I15").Value = 0.4664013802790521
$ws.Range("J15").Value = 0.0321582140825063
$ws.Range("M15").Value = 0.672341843352612
$ws.Range("N15").Value = 1.04201792785091
$ws.Range("B16").Value = 1.037874408996458
$ws.Range("C16").Value = 0.292967574422704
$ws.Range("E16").Value = 0.4731990280769622
$ws.Range("F16").Value = 2.178923892475012
$ws.Range("G16").Value = 0.6118891200861754
$ws.Range("H16").Value = 0.6596865306965043
$ws.Range("I16").Value = 0.4666874900090079
$ws.Range("J16").Value = 0.03251939315157948
$ws.Range("M16").Value = 0.6482539749360967
$ws.Range("N16").Value = 1.050841105706787
$ws.Range("B17").Value = 0.997743015256674
$ws.Range("C17").Value = 0.2810652583083311
$ws.Range("E17").Value = 0.4716460847558821
$ws.Range("F17").Value = 2.166159164025302
$ws.Range("G17").Value = 0.6054348995154584
$ws.Range("H17").Value = 0.6586932124970701
$ws.Range("I17").Value = 0.4669653968835732
$ws.Range("J17").Value = 0.03274657290599281
$ws.Range("M17").Value = 0.6335339035435368
$ws.Range("N17").Value = 1.056387098978263
$ws.Range("B18").Value = 0.9746628683333824
$ws.Range("C18").Value = 0.2742166953857463
$ws.Range("E18").Value = 0.4707689959410217
$ws.Range("F18").Value = 2.158945071866881
$ws.Range("G18").Value = 0.6017769896696308
$ws.Range("H18").Value = 0.6581648629972392
$ws.Range("I18").Value = 0.4671627625335617
$ws.Range("J18").Value = 0.03287929751995833
$ws.Range("M18").Value = 0.6250883398561768
$ws.Range("N18").Value = 1.059625931617624
$ws.Range("B19").Value = 0.96684875134099
$ws.Range("C19").Value = 0.2718974315822891
$ws.Range("E19").Value = 0.4704747999075778
$ws.Range("F19").Value = 2.156524444266452
$ws.Range("G19").Value = 0.6005478025804507
$ws.Range("H19").Value = 0.657993347462039
$ws.Range("I19").Value = 0.4672360213563849
$ws.Range("J19").Value = 0.03292458866615444
$ws.Range("M19").Value = 0.6222324355045146
$ws.Range("N19").Value = 1.060730942920326
$ws.Range("B20").Value = 1.002014829512632
$ws.Range("C20").Value = 0.2823325559845387
$ws.Range("E20").Value = 0.4718097300192028
$ws.Range("F20").Value = 2.167504755949565
$ws.Range("G20").Value = 0.6061163280488131
$ws.Range("H20").Value = 0.6587945027736168
$ws.Range("I20").Value = 0.4669319275477335
$ws.Range("J20").Value = 0.03272217626639984
$ws.Range("M20").Value = 0.6350987049893746
$ws.Range("N20").Value = 1.055791653991019
$ws.Range("B21").Value = 1.120195073334969
$ws.Range("C21").Value = 0.3173613192815594
$ws.Range("E21").Value = 0.4764874005158717
$ws.Range("F21").Value = 2.205925792302708
$ws.Range("G21").Value = 0.6254773005959322
$ws.Range("H21").Value = 0.6620003578208298
$ws.Range("I21").Value = 0.4663592533349075
$ws.Range("J21").Value = 0.03206670192560468
$ws.Range("M21").Value = 0.6785784447056074
$ws.Range("N21").Value = 1.03978119529642
$ws.Range("B22").Value = 1.197449958001243
$ws.Range("C22").Value = 0.3402309736580378
$ws.Range("E22").Value = 0.4796876568279203
$ws.Range("F22").Value = 2.23217840106318
$ws.Range("G22").Value = 0.6386202226570106
$ws.Range("H22").Value = 0.6644799934041998
$ws.Range("I22").Value = 0.466321356704988
$ws.Range("J22").Value = 0.03165671537366599
$ws.Range("M22").Value = 0.7071808965346946
$ws.Range("N22").Value = 1.029754138680943
$ws.Range("B23").Value = 1.156216176973885
$ws.Range("C23").Value = 0.3280271706014162
$ws.Range("E23").Value = 0.4779665450012658
$ws.Range("F23").Value = 2.218062188542376
$ws.Range("G23").Value = 0.6315605698861333
$ws.Range("H23").Value = 0.663121308511009
$ws.Range("I23").Value = 0.4663106945363182
$ws.Range("J23").Value = 0.03187384467484033
$ws.Range("M23").Value = 0.6918982825354334
$ws.Range("N23").Value = 1.035065757613506
$ws.Range("B24").Value = 1.000083567398633
$ws.Range("C24").Value = 0.281759628662428
$ws.Range("E24").Value = 0.4717356970274693
$ws.Range("F24").Value = 2.166896026001538
$ws.Range("G24").Value = 0.6058080900829452
$ws.Range("H24").Value = 0.6587485763735117
$ws.Range("I24").Value = 0.4669469419542267
$ws.Range("J24").Value = 0.03273319939676078
$ws.Range("M24").Value = 0.6343912046613553
$ws.Range("N24").Value = 1.0560606977809
$ws.Range("B25").Value = 0.8317615546712886
$ws.Range("C25").Value = 0.2317504768931826
$ws.Range("E25").Value = 0.4656313728744621
$ws.Range("F25").Value = 2.116590172454536
$ws.Range("G25").Value = 0.5801060481821736
$ws.Range("H25").Value = 0.6556740891696364
$ws.Range("I25").Value = 0.4690663846528906
$ws.Range("J25").Value = 0.03373874022975176
$ws.Range("M25").Value = 0.5731651910737767
$ws.Range("N25").Value = 1.08057708726961
